# Fix Training Data Issue (#48)
# The "Date" column (BF) on Sheet1 held values like "6-12-2013-14" which
# mixes the game date with the season label. Re-write it as a clean
# ISO-ish date string "2014-06-12" for every data row (BF2:BF31).
#
# A plain Range.Value assignment of "2014-06-12" gets auto-recognized by
# Excel as a real date and silently converted to a date serial number,
# which is not what we want (the source file stores these as literal
# text). To keep the literal text without permanently changing the
# number format of the target cells, stage the text in a scratch cell
# that is explicitly formatted as Text, copy it, and paste-special just
# the values into the target cells (paste values only does not carry the
# source number format along), then discard the scratch column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 31
$dateColumn = "BF"
$oldValue = "6-12-2013-14"
$newValue = "2014-06-12"

$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = $newValue
$scratch.Copy()

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $cell = $ws.Range("$dateColumn$row")
    if ($cell.Text -eq $oldValue) {
        $cell.PasteSpecial(-4163)
    }
}

$excel.CutCopyMode = $false
$scratch.EntireColumn.Delete()
